# "Added one dot to image" - split the title run on slide 1 so the
# trailing "X" run becomes its own run with an extra "X" added in
# front of it: "GIT vs X…." -> "GIT vs " + "XX…."

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$titleShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -like "GIT vs*") {
            $titleShape = $sh
            break
        }
    }
}

$tr = $titleShape.TextFrame.TextRange
$firstRun = $tr.Runs(1)

# Keep "GIT vs " on the original run, and add the remainder ("XX….")
# as a brand-new run straight after it.
$firstRun.Text = "GIT vs "
$tr.InsertAfter("XX….") | Out-Null
